# Update on 23/03/2020 at 15:18
#
# The document has a stray paragraph (right before the "_GoBack" bookmark)
# that only contains the placeholder run "test". Remove that run's text
# so the paragraph is left empty (its paragraph mark, pPr/rPr formatting,
# and the _GoBack bookmarkStart/bookmarkEnd all stay in place) -- exactly
# mirroring a user selecting the word and pressing Delete.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.Text = "test"
$find.Replacement.Text = ""
$find.Execute(
    "test",   # FindText
    $true,    # MatchCase
    $false,   # MatchWholeWord
    $false,   # MatchWildcards
    $false,   # MatchSoundsLike
    $false,   # MatchAllWordForms
    $true,    # Forward
    1,        # Wrap (wdFindContinue)
    $false,   # Format
    "",       # ReplaceWith
    2         # Replace (wdReplaceAll)
)
